$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.022961123942673686
$ws.Range("B3").Value = 0.13100861321262741
$ws.Range("B4").Value = 0.10966941863989707
$ws.Range("B5").Value = 0.031739309075173007
$ws.Range("B6").Value = 0.24697397787159978
$ws.Range("B7").Value = 0.12417711935984306
$ws.Range("B8").Value = 0.019206265986033866
$ws.Range("B9").Value = 0.36450065340845694
